$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) First, append the two new rows (160,161) as an exact copy of the
#    current rows 158:159 (before any of the date/origin values below are
#    changed). These rows keep their original date (44425) and origin
#    (Region de Nuble), matching the new tail of the shifted sequence.
$ws.Range("A158:R159").Copy($ws.Range("A160:R161"))

# 2) Update the "Fecha" (column D) values for rows 138-159. The weekly
#    date/origin sequence shifts down by one pair and a new pair is
#    inserted at the top (rows 138:139, date 44474).
$ws.Range("D138").Value = 44474
$ws.Range("D139").Value = 44474
$ws.Range("D140").Value = 44209
$ws.Range("D141").Value = 44209
$ws.Range("D142").Value = 44365
$ws.Range("D143").Value = 44365
$ws.Range("D144").Value = 44306
$ws.Range("D145").Value = 44306
$ws.Range("D146").Value = 44215
$ws.Range("D147").Value = 44215
$ws.Range("D148").Value = 44257
$ws.Range("D149").Value = 44257
$ws.Range("D150").Value = 44239
$ws.Range("D151").Value = 44239
$ws.Range("D152").Value = 44376
$ws.Range("D153").Value = 44376
$ws.Range("D154").Value = 44292
$ws.Range("D155").Value = 44292
$ws.Range("D156").Value = 44358
$ws.Range("D157").Value = 44358
$ws.Range("D158").Value = 44211
$ws.Range("D159").Value = 44211

# 3) Update the "Origen" (column O) values that change as part of the
#    same shift (rows 148:149 and 150:151 swap region).
$ws.Range("O148").Value = "Región de Ñuble"
$ws.Range("O149").Value = "Región de Ñuble"
$ws.Range("O150").Value = "Región Metropolitana"
$ws.Range("O151").Value = "Región Metropolitana"
